$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-6 (Generation 0-4): Fitness column C changes from 7651 to 7295
$ws.Range("C2:C6").Value = 7295

# Rows 7-252 (Generation 5-250): Fitness column C changes from 7651 to 7293
$ws.Range("C7:C252").Value = 7293
